# Update generated Excel file
#
# Sheet "autonomous_t26" (the workbook's first/active sheet) gains one more
# CAN message block at the bottom - "Message: dir" (ID 0x1, Sender(s):
# Unknown) with a single signal "dir" - formatted exactly like every other
# message block already on the sheet (blue message-header row, yellow
# column-header row, bordered data row), separated from the previous block
# by one blank spacer row. Column C also widens a bit (19 -> 20) to fit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- widen column C (Length (bits)) from 19 to 20 -------------------------
$ws.Columns.Item(3).ColumnWidth = 19.2

# --- new message block, rows 24-26 (row 23 is left blank, like the other
#     spacer rows already on the sheet, e.g. row 11 / row 19) -------------

# Pull down the formatting used by the previous block (rows 20-22) so the
# new block's styling (fills/borders/bold) matches the rest of the sheet.
$ws.Range("A20:C20").Copy($ws.Range("A24:C24"))
$ws.Range("A21:K21").Copy($ws.Range("A25:K25"))
$ws.Range("A22:K22").Copy($ws.Range("A26:K26"))

# Message header row (24).
$ws.Range("A24").Value = "Message: dir"
$ws.Range("B24").Value = "ID: 0x1"
$ws.Range("C24").Value = "Sender(s): Unknown"

# Column header row (25) already reads Signal Name / Start Bit / ... /
# Choices after the copy from row 21 - nothing else to change there.

# Signal data row (26) - signal "dir".
$ws.Range("A26").Value = "dir"
$ws.Range("B26").Value = 0
$ws.Range("C26").Value = 8
$ws.Range("D26").Value = "Intel"
$ws.Range("E26").Value = $false
$ws.Range("F26").Value = 1
$ws.Range("G26").Value = 0
$ws.Range("H26").ClearContents()
$ws.Range("I26").ClearContents()
$ws.Range("J26").ClearContents()
$ws.Range("K26").ClearContents()
